$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.700.68"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.531.30"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.43"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.36"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0579"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "1.749.34"
$ws.Range("D13").Value = "1.539.93"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "26.691.50"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.03"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.05"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "0.0₃0682"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.75"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("D33").Value = "1.359.98"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.938"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0163"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.519"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.70"
$ws.Range("E41").Value = "  +6.14%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.797"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.20"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "1.663.60"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.40"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0504"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -0.76%  "
